# Update 2021 Target Depth Data: log divisional round, simulating season from conference round.
$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 248
$wsOff.Range("C3").Value = 172
$wsOff.Range("D3").Value = 65
$wsOff.Range("E3").Value = 38
$wsOff.Range("F3").Value = 5

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 254
$wsDef.Range("C3").Value = 189
$wsDef.Range("D3").Value = 57
$wsDef.Range("E3").Value = 27
